# Clarify status report for unit tests: distinguish an "ERROR" result
# (the expected/actual comparison itself errored out) from a plain "FAIL".
#
# D3 holds its own (non-shared) formula; D4:D29 share formula group si="0".
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D3").Formula = '=IF(ISERROR(B3),"ERROR",IF(ISERROR(C3),"FAIL",IF(B3=C3,"PASS","FAIL")))'
$ws.Range("D4:D29").Formula = '=IF(ISERROR(B4),"ERROR",IF(ISERROR(C4),"FAIL",IF(B4=C4,"PASS","FAIL")))'
